# Updates crypto price/volume figures (and two pairs of swapped coin rows)
# to match the "Updated cryptos list" GitHub Actions commit.
# Note: some new price strings (e.g. "616.32") look like plain numbers, so a
# leading apostrophe is used to force Excel to keep them as text, matching
# the original inlineStr cell type instead of letting Excel coerce them to
# numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.903.57'
$ws.Range('E2').Value = '  +1.04%  '
$ws.Range('D3').Value = '3.542.39'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '''616.32'
$ws.Range('E5').Value = '  +0.54%  '
$ws.Range('D6').Value = '''152.94'
$ws.Range('E6').Value = '  -1.23%  '
$ws.Range('D7').Value = '3.542.75'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = '''0.482'
$ws.Range('E9').Value = '  -0.89%  '
$ws.Range('D10').Value = '''0.141'
$ws.Range('E10').Value = '  -1.06%  '
$ws.Range('D11').Value = '''7.10'
$ws.Range('E11').Value = '  +3.30%  '
$ws.Range('E12').Value = '  -0.99%  '
$ws.Range('E13').Value = '  -0.72%  '
$ws.Range('D14').Value = '4.141.81'
$ws.Range('E14').Value = '  -0.03%  '
$ws.Range('D15').Value = '''32.18'
$ws.Range('E15').Value = '  +0.13%  '
$ws.Range('D16').Value = '3.554.63'
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('D17').Value = '67.660.90'
$ws.Range('E17').Value = '  +0.67%  '
$ws.Range('D19').Value = '''6.42'
$ws.Range('E19').Value = '  +0.17%  '
$ws.Range('D20').Value = '''15.39'
$ws.Range('E20').Value = '  -0.52%  '
$ws.Range('D21').Value = '''448.02'
$ws.Range('E21').Value = '  -1.16%  '
$ws.Range('D22').Value = '''9.68'
$ws.Range('E22').Value = '  +2.54%  '
$ws.Range('D23').Value = '''0.626'
$ws.Range('E23').Value = '  -2.86%  '
$ws.Range('D24').Value = '''77.61'
$ws.Range('E24').Value = '  -2.34%  '
$ws.Range('E25').Value = '  +6.47%  '
$ws.Range('D26').Value = '3.684.02'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('D28').Value = '''10.30'
$ws.Range('E28').Value = '  -1.23%  '
$ws.Range('D29').Value = '''8.69'
$ws.Range('E29').Value = '  +3.61%  '
$ws.Range('D30').Value = '''2.54'
$ws.Range('E30').Value = '  -1.21%  '
$ws.Range('E31').Value = '  -4.22%  '
$ws.Range('E32').Value = '  +7.11%  '
$ws.Range('D33').Value = '''0.999'
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('D34').Value = '''25.99'
$ws.Range('E34').Value = '  -0.10%  '
$ws.Range('D35').Value = '''6.23'
$ws.Range('E35').Value = '  +0.29%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '''1.86'
$ws.Range('E36').Value = '  -2.25%  '
$ws.Range('B37').Value = 'RenzoRestakedETH'
$ws.Range('C37').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D37').Value = '3.527.29'
$ws.Range('E37').Value = '  -0.36%  '
$ws.Range('D38').Value = '''8.07'
$ws.Range('E38').Value = '  -0.30%  '
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').Value = '''0.999'
$ws.Range('E40').Value = '  -0.11%  '
$ws.Range('D41').Value = '''176.68'
$ws.Range('E41').Value = '  -0.48%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '''2.19'
$ws.Range('E42').Value = '  +2.26%  '
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').Value = '''0.0896'
$ws.Range('E43').Value = '  +2.05%  '
$ws.Range('D44').Value = '''5.44'
$ws.Range('E44').Value = '  -3.63%  '
$ws.Range('E45').Value = '  -0.82%  '
$ws.Range('D46').Value = '''28.59'
$ws.Range('E46').Value = '  -0.61%  '
$ws.Range('E47').Value = '  -1.05%  '
$ws.Range('D48').Value = '''2.69'
$ws.Range('E48').Value = '  -0.54%  '
$ws.Range('D49').Value = '''1.29'
$ws.Range('E49').Value = '  +5.32%  '
$ws.Range('D50').Value = '''7.64'
$ws.Range('E50').Value = '  -0.64%  '
$ws.Range('D51').Value = '''1.00'
$ws.Range('E51').Value = '  -3.35%  '
